$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: servings 41 -> 40
$ws.Range("C2").Value = 40

# Row 3: servings 29 -> 100
$ws.Range("C3").Value = 100

# Row 6: id 6 -> 7, stock_type "injeer" -> "Creatine ", servings 66 -> 93,
# cost_per_serving 55 -> 150, profit_per_serving 45 -> 50
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Creatine "
$ws.Range("C6").Value = 93
$ws.Range("D6").Value = 150
$ws.Range("E6").Value = 50

# Row 7: stock_type "Creatine " -> "protine powder", servings 93 -> 150,
# cost_per_serving 150 -> 100, profit_per_serving 50 -> 20,
# other_charges 10 -> 5, date_added 2025-03-10 -> 2025-03-13
$ws.Range("B7").Value = "protine powder"
$ws.Range("C7").Value = 150
$ws.Range("D7").Value = 100
$ws.Range("E7").Value = 20
$ws.Range("F7").Value = 5
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "2025-03-13"
$ws.Range("G6").Copy()
$ws.Range("G7").PasteSpecial(-4122)  # xlPasteFormats, restore plain style (no date auto-format)
$excel.CutCopyMode = 0
